# Parameters sheet: insert two new blank columns (M:N) ahead of the
# existing "Population type" column (which shifts from M to O), and
# label the two new columns' headers "Unnamed: 12" / "Unnamed: 13".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

$ws.Columns("M:N").Insert()

$ws.Range("M1").Value = "Unnamed: 12"
$ws.Range("N1").Value = "Unnamed: 13"
